# Merge the two paragraphs:
#   "A new chapter begins."
#   "The story continues with new developments and revelations. Characters
#    face challenges and make important decisions that will shape the
#    narrative going forward."
# into a single paragraph/run:
#   "A new chapter begins. The story continues with new developments and
#    revelations. Characters face challenges and make important decisions
#    that will shape the narrative going forward."

$d = $word.ActiveDocument

$findText = "A new chapter begins.^p" + `
    "The story continues with new developments and revelations. " + `
    "Characters face challenges and make important decisions that will " + `
    "shape the narrative going forward."

$replaceText = "A new chapter begins. " + `
    "The story continues with new developments and revelations. " + `
    "Characters face challenges and make important decisions that will " + `
    "shape the narrative going forward."

$range = $d.Content
$found = $range.Find.Execute($findText, $true, $false, $false, $false, $false, `
    $true, 1, $false, $replaceText, 2)

if (-not $found) {
    throw "Could not find the target paragraphs to merge."
}

Write-Host "Merged paragraph found and replaced: $found"
